$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 77, column B: change from text "3" to a real number 3
$ws.Cells.Item(77, 2).Value = 3

# Add new row 78 with data
$ws.Cells.Item(78, 1).Value = "Ruilin"
$ws.Cells.Item(78, 2).Value = "'3"
$ws.Cells.Item(78, 2).Style = "Normal"
$ws.Cells.Item(78, 3).Value = "无"
$ws.Cells.Item(78, 4).Value = "DIS"
$ws.Cells.Item(78, 5).Value = "RES"
$ws.Cells.Item(78, 6).Value = "0a719846-8538-4eb8-b511-cc6f0b597898"
$ws.Cells.Item(78, 7).Value = "bb7SwHahSUpiq_annotated.xlsx"
$ws.Cells.Item(78, 8).Value = "Fig. 1 (a)-(d) agree as well with my own experience of training GRBMs on image patches."
